$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as text in this sheet.
# Force text format first so Excel does not auto-convert numeric-looking
# strings (e.g. "581.58", "1.00", "0.0407") into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.625.11'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '3.337.41'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '581.58'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').Value = '176.38'
$ws.Range('E6').Value = '  -3.50%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '3.334.46'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').Value = '0.578'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '45.51'
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('E13').Value = '  -2.50%  '
$ws.Range('D14').Value = '665.30'
$ws.Range('E14').Value = '  +4.09%  '
$ws.Range('D15').Value = '3.872.73'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '8.43'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '67.706.06'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '3.334.01'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '17.43'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').Value = '10.99'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = '0.891'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = '5.50'
$ws.Range('E23').Value = '  +9.30%  '
$ws.Range('D24').Value = '17.10'
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('D25').Value = '99.61'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').Value = '  -3.57%  '
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -4.80%  '
$ws.Range('D28').Value = '9.32'
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('D29').Value = '33.72'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').Value = '7.45'
$ws.Range('E30').Value = '  +11.30%  '
$ws.Range('D31').Value = '8.46'
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').Value = '578.38'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').Value = '10.99'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').Value = '0.105'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').Value = '3.710.56'
$ws.Range('E36').Value = '  -5.76%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = '56.64'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '3.39'
$ws.Range('E38').Value = '  -6.46%  '
$ws.Range('D39').Value = '34.51'
$ws.Range('E39').Value = '  +4.33%  '
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('D41').Value = '2.63'
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('D42').Value = '3.12'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('E43').Value = '  -2.64%  '
$ws.Range('D44').Value = '0.335'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('D45').Value = '3.29'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').Value = '0.0407'
$ws.Range('E46').Value = '  -2.18%  '
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').Value = '128.35'
$ws.Range('E51').Value = '  -2.02%  '
